# Issue #57: add a required "Genre" column (PBCore controlled vocabulary)
# to the batch-manifest spreadsheet fixture, alongside the existing
# "Topical Subject" column, with a sample value of "Auction".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header in row 2 (row 1 holds the submitter name/email), right after
# the existing "Topical Subject" header in column I.
$ws.Range("J2").Value = "Genre"

# Sample data row (row 3) value for the new Genre column.
$ws.Range("J3").Value = "Auction"

# Match the author's final selection/cursor position in the saved file.
$ws.Range("J3").Select()
